# Add a few more java practice programs to the tracking sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in ProgramName (B) and Status (C) for rows 4-7, which previously
# only had the "Programs" description filled in column A.
$ws.Range("B4").Value = "multiply2Floating.java"
$ws.Range("C4").Value = "DONE"

$ws.Range("B5").Value = "findASCII.java"
$ws.Range("C5").Value = "DONE"

$ws.Range("B6").Value = "findQuotientAndRemainder.java"
$ws.Range("C6").Value = "DONE"

$ws.Range("B7").Value = "swapUsing3Variable.java"
$ws.Range("C7").Value = "DONE"

# Update the view so the selection ends up on C7, with the sheet scrolled
# down so row 4 is the first visible row.
$ws.Range("C7").Select()
$excel.ActiveWindow.ScrollRow = 4
